# Update "想去人数" (want-to-go count) figures in the F column
# for the sheets that contain the event data ("展览" and "全部类型").
# Both sheets hold identical data, so the same cell updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 2663
    "F7"  = 139
    "F9"  = 1403
    "F13" = 1200
    "F22" = 2566
    "F23" = 36
    "F24" = 293
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
